# Customers sheet: refresh customer list (signin/signout backend data update)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Customers")

# Remove the mailto hyperlink currently anchored on D2 (email column) -
# the refreshed data no longer carries a live hyperlink on that cell.
$ws.Range("D2").Hyperlinks.Delete()

# Column E (Password) and Column C (Mobile) contain numeric-looking text
# ("1234", phone numbers) that must be stored as text, not numbers, so
# force a text number format on the whole data block before writing values.
$dataRange = $ws.Range("A2:G5")
$dataRange.NumberFormat = "@"

# Row 2: C#00003 / Kishore kumar (moved up from the old row 4)
$ws.Range("A2").Value = "C#00003"
$ws.Range("B2").Value = "Kishore kumar"
$ws.Range("C2").Value = "0123789456"
$ws.Range("D2").Value = "kishorekumar@gmail.com"
$ws.Range("E2").Value = "1234"
$ws.Range("F2").Value = "Patna"
$ws.Range("G2").Value = "ACTIVE"

# Row 3: new customer C#00004 / Aman
$ws.Range("A3").Value = "C#00004"
$ws.Range("B3").Value = "Aman"
$ws.Range("C3").Value = "9935465163"
$ws.Range("D3").Value = "amankumar@gmail.com"
$ws.Range("E3").Value = "1234"
$ws.Range("F3").Value = "Patna"
$ws.Range("G3").Value = "ACTIVE"

# Row 4: new customer C#00005 / Paramjot Singh
$ws.Range("A4").Value = "C#00005"
$ws.Range("B4").Value = "Paramjot Singh"
$ws.Range("C4").Value = "9031398069"
$ws.Range("D4").Value = "paramjotsingh966@gmail.com"
$ws.Range("E4").Value = "1234"
$ws.Range("F4").Value = "Pragathi Layout"
$ws.Range("G4").Value = "ACTIVE"

# Row 5: new customer C#00006 / Abhi
$ws.Range("A5").Value = "C#00006"
$ws.Range("B5").Value = "Abhi"
$ws.Range("C5").Value = "1234567890"
$ws.Range("D5").Value = "abhi@gmail.com"
$ws.Range("E5").Value = "1234"
$ws.Range("F5").Value = "Marathahalli"
$ws.Range("G5").Value = "ACTIVE"

# Strip the now-unneeded cell styles (bold/hyperlink/text-format xf indexes)
# from the data rows so plain default-styled cells remain, matching the
# refreshed sheet's formatting.
$dataRange.ClearFormats()
